$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two missing values in column D (rows 9 and 10) to match
# the rest of the table's data.
$ws.Range("D9").Value = 0.0
$ws.Range("D10").Value = 0.0

# Right-align the main data block (A1:F12) — previously a mix of
# centered and default alignment.
$ws.Range("A1:F12").HorizontalAlignment = -4152
